# Corrected excel sheets for application fix issues
#
# This script reproduces, via Excel COM-interop, the edits that were made to
# the "Makerepayment1" workbook:
#   - The no-longer-needed "Acc_Disbursement" and "Acc_Repayment" ledger
#     sheets are removed.
#   - The repayment amounts (which were wrong - interest/fee split, totals,
#     running balances) are corrected on "Summary", "Repayment schedule" and
#     "Transactions".
#   - The "Transactions" sheet trims its stray extra/blank row.
#   - Selections on a couple of sheets are moved to reflect where the user
#     was last working.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# 1. Remove the two accounting sheets that are no longer part of the test
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Acc_Disbursement").Delete()
$wb.Worksheets.Item("Acc_Repayment").Delete()

# ---------------------------------------------------------------------
# 2. Summary sheet - fix the "Fees" row (row 5) amounts
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A5").Value = 0.89
$wsSummary.Range("B5").Value = 0.51
$wsSummary.Range("E5").Value = 0.38
$wsSummary.Range("F5").Value = 0.38
$wsSummary.Range("C5").Select()

# ---------------------------------------------------------------------
# 3. Repayment schedule - correct the fee/interest split, the due/paid
#    totals and the running "Outstanding" balance for the first
#    repayment installments
# ---------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Row 2 (opening/disbursement row): the trailing blank "Outstanding" cell
# shifts from column P to column O
$wsSchedule.Cells.Item(2, 14).Copy()
$wsSchedule.Cells.Item(2, 15).PasteSpecial(-4122)   # xlPasteFormats
$wsSchedule.Cells.Item(2, 16).Clear()               # P2 removed

# Row 3 (installment 1)
$wsSchedule.Cells.Item(3, 10).Value = 0.51      # J3 Fees Due
$wsSchedule.Cells.Item(3, 11).Value = 888.23    # K3 Due
$wsSchedule.Cells.Item(3, 12).Value = 888.23    # L3 Paid
$wsSchedule.Cells.Item(3, 15).Clear()           # O3 (extra outstanding-fee column removed)

# Row 4 (installment 2)
$wsSchedule.Cells.Item(4, 10).Value = 0.38      # J4 Fees Due
$wsSchedule.Cells.Item(4, 11).Value = 888.1     # K4 Due
$wsSchedule.Cells.Item(4, 12).Value = 0         # L4 Paid
$wsSchedule.Cells.Item(4, 15).Clear()           # O4
$wsSchedule.Cells.Item(4, 16).Value = 888.1     # P4 Outstanding

# Row 5 (installment 3)
$wsSchedule.Cells.Item(5, 10).Value = 0         # J5 Fees Due
$wsSchedule.Cells.Item(5, 11).Value = 887.72    # K5 Due
$wsSchedule.Cells.Item(5, 15).Clear()           # O5
$wsSchedule.Cells.Item(5, 16).Value = 887.72    # P5 Outstanding

# Rows 6-8: the now-unused "O" (outstanding fee) column is cleared
$wsSchedule.Cells.Item(6, 15).Clear()           # O6
$wsSchedule.Cells.Item(7, 15).Clear()           # O7
$wsSchedule.Cells.Item(8, 15).Clear()           # O8

$wsSchedule.Range("F8").Select()

# ---------------------------------------------------------------------
# 4. Input sheet - no longer the active tab, but selection stays on A2
# ---------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Activate()
$wsInput.Range("A2").Select()

# ---------------------------------------------------------------------
# 5. Transactions sheet - correct the repayment transaction row and drop
#    the stray blank row that trailed the data. Ends up the active tab.
# ---------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")

$wsTransactions.Cells.Item(2, 1).Value = 1738          # A2 Entry Id
$wsTransactions.Cells.Item(2, 4).Value = "Repayment"   # D2 Transaction Type
$wsTransactions.Cells.Item(2, 5).Value = 888.23        # E2 Amount
$wsTransactions.Cells.Item(2, 9).Value = 0.51          # I2 Fees
$wsTransactions.Cells.Item(2, 10).Value = 4163.24      # J2 Loan Balance

$wsTransactions.Cells.Item(3, 1).Value = 678           # A3 Entry Id

# Remove the stray extra row 6 (J6) left over below the data
$wsTransactions.Rows.Item(6).Delete()

$wsTransactions.Columns.Item(1).ColumnWidth = 5

$wsTransactions.Activate()
$wsTransactions.Range("F3").Select()

$wb.Save()
